# DT-1853 update template to reflect new field descriptions on instructions
#
# The "project_collection_template.xlsx" had its generic "(string)"
# placeholder text replaced with field-specific placeholder text on both
# the "Data collection project" and "Data collection" sheets (Title,
# Short title and Description rows), and the active sheet/selection
# moved from "Data collection project" to "Data collection".

$wb = $excel.ActiveWorkbook

$wsProject    = $wb.Worksheets.Item("Data collection project")
$wsCollection = $wb.Worksheets.Item("Data collection")

# --- "Data collection project" sheet: field placeholders ---------------
$wsProject.Range("B13").Value = "(project title)"
$wsProject.Range("B14").Value = "(project short title)"
$wsProject.Range("B15").Value = "(project description)"

# --- "Data collection" sheet: field placeholders ------------------------
$wsCollection.Range("B8").Value  = "(collection title)"
$wsCollection.Range("B9").Value  = "(collection short title)"
$wsCollection.Range("B10").Value = "(collection description)"

# --- View-state: column B a touch wider on "Data collection" -----------
$wsCollection.Columns.Item(2).ColumnWidth = 20.42

# --- View-state: selection on "Data collection project" ----------------
$wsProject.Activate() | Out-Null
$wsProject.Range("C11").Select() | Out-Null

# --- View-state: "Data collection" becomes the active tab --------------
$wsCollection.Activate() | Out-Null
$wsCollection.Range("D8").Select() | Out-Null
